# Updates the crypto price/volume table (columns D and E) to the latest
# scraped values, mirroring the GitHub Actions data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new Price (column D, or $null if unchanged),
# new Volume(1h) percentage text without padding (column E, or $null if unchanged).
$updates = @(
    [PSCustomObject]@{ Row = 2; Price = "64.061.58"; Volume = "-3.55%" }
    [PSCustomObject]@{ Row = 3; Price = "3.161.28"; Volume = "-8.75%" }
    [PSCustomObject]@{ Row = 5; Price = "562.58"; Volume = "-4.05%" }
    [PSCustomObject]@{ Row = 6; Price = "167.98"; Volume = "-4.89%" }
    [PSCustomObject]@{ Row = 7; Price = $null; Volume = "+0.01%" }
    [PSCustomObject]@{ Row = 8; Price = "0.604"; Volume = "-3.36%" }
    [PSCustomObject]@{ Row = 9; Price = "3.159.21"; Volume = "-8.71%" }
    [PSCustomObject]@{ Row = 10; Price = $null; Volume = "-7.14%" }
    [PSCustomObject]@{ Row = 11; Price = $null; Volume = "-5.25%" }
    [PSCustomObject]@{ Row = 12; Price = $null; Volume = "-5.63%" }
    [PSCustomObject]@{ Row = 13; Price = "3.705.73"; Volume = "-8.84%" }
    [PSCustomObject]@{ Row = 14; Price = $null; Volume = "+1.40%" }
    [PSCustomObject]@{ Row = 15; Price = "27.23"; Volume = "-10.28%" }
    [PSCustomObject]@{ Row = 16; Price = "64.052.42"; Volume = "-3.40%" }
    [PSCustomObject]@{ Row = 17; Price = $null; Volume = "-5.66%" }
    [PSCustomObject]@{ Row = 18; Price = "3.160.00"; Volume = "-8.96%" }
    [PSCustomObject]@{ Row = 19; Price = "5.71"; Volume = "-4.46%" }
    [PSCustomObject]@{ Row = 20; Price = "12.93"; Volume = "-6.41%" }
    [PSCustomObject]@{ Row = 21; Price = "351.00"; Volume = "-5.93%" }
    [PSCustomObject]@{ Row = 22; Price = $null; Volume = "-6.69%" }
    [PSCustomObject]@{ Row = 23; Price = "0.999"; Volume = "+0.00%" }
    [PSCustomObject]@{ Row = 24; Price = "68.54"; Volume = "-6.46%" }
    [PSCustomObject]@{ Row = 25; Price = "0.0000119"; Volume = "-6.15%" }
    [PSCustomObject]@{ Row = 26; Price = $null; Volume = "-6.63%" }
    [PSCustomObject]@{ Row = 27; Price = "9.50"; Volume = "-4.43%" }
    [PSCustomObject]@{ Row = 28; Price = $null; Volume = "-1.55%" }
    [PSCustomObject]@{ Row = 29; Price = $null; Volume = "-0.07%" }
    [PSCustomObject]@{ Row = 30; Price = $null; Volume = "-0.16%" }
    [PSCustomObject]@{ Row = 31; Price = $null; Volume = "-7.04%" }
    [PSCustomObject]@{ Row = 32; Price = $null; Volume = "-5.72%" }
    [PSCustomObject]@{ Row = 33; Price = "21.90"; Volume = $null }
    [PSCustomObject]@{ Row = 34; Price = $null; Volume = "-6.31%" }
    [PSCustomObject]@{ Row = 35; Price = $null; Volume = "-5.71%" }
    [PSCustomObject]@{ Row = 36; Price = $null; Volume = "-7.74%" }
    [PSCustomObject]@{ Row = 37; Price = "153.49"; Volume = "-4.77%" }
    [PSCustomObject]@{ Row = 38; Price = "0.812"; Volume = "-8.32%" }
    [PSCustomObject]@{ Row = 39; Price = "25.63"; Volume = "-9.34%" }
    [PSCustomObject]@{ Row = 40; Price = "2.50"; Volume = "-3.18%" }
    [PSCustomObject]@{ Row = 41; Price = $null; Volume = "-6.63%" }
    [PSCustomObject]@{ Row = 42; Price = "2.595.22"; Volume = "-6.81%" }
    [PSCustomObject]@{ Row = 43; Price = $null; Volume = "-7.68%" }
    [PSCustomObject]@{ Row = 44; Price = "39.33"; Volume = "-1.76%" }
    [PSCustomObject]@{ Row = 45; Price = "5.97"; Volume = "-7.56%" }
    [PSCustomObject]@{ Row = 46; Price = $null; Volume = "-6.52%" }
    [PSCustomObject]@{ Row = 47; Price = "23.62"; Volume = "-6.74%" }
    [PSCustomObject]@{ Row = 48; Price = "317.31"; Volume = "-6.76%" }
    [PSCustomObject]@{ Row = 49; Price = $null; Volume = "-8.60%" }
    [PSCustomObject]@{ Row = 50; Price = $null; Volume = "-3.39%" }
    [PSCustomObject]@{ Row = 51; Price = $null; Volume = "-0.01%" }
)

foreach ($u in $updates) {
    if ($null -ne $u.Price) {
        $priceCell = $ws.Range("D" + $u.Row)
        # Force text so multi-group numbers (e.g. 64.061.58) and plain
        # decimals (e.g. 0.604) are both stored as text, matching the
        # original inline-string cell content instead of being parsed
        # into a numeric value.
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $u.Price
        $priceCell.Style = "Normal"
    }
    if ($null -ne $u.Volume) {
        $volumeCell = $ws.Range("E" + $u.Row)
        $volumeCell.Value = "  " + $u.Volume + "  "
    }
}
